$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (Volume(1h) %) updates ---
# These are padded text values (e.g. "  -0.89%  ") so Excel keeps them as
# plain text automatically - no special handling required.
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("E6").Value = "  -4.24%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("E10").Value = "  -6.32%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  -2.18%  "
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").Value = "  -8.74%  "
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("E40").Value = "  -5.46%  "
$ws.Range("E41").Value = "  -5.74%  "
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  -4.09%  "
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("E51").Value = "  -0.84%  "

# --- Column D (Price) updates ---
# Price text like "72.170.80" (grouped thousands, not a valid Excel number)
# is safe to assign directly; it stays text.
$ws.Range("D2").Value = "72.170.80"
$ws.Range("D3").Value = "2.691.56"
$ws.Range("D9").Value = "2.690.12"
$ws.Range("D14").Value = "3.183.42"
$ws.Range("D16").Value = "72.070.85"
$ws.Range("D18").Value = "2.685.93"
$ws.Range("D28").Value = "2.830.03"

# These new price strings (e.g. "599.23") parse as valid numbers, so a plain
# .Value assignment would silently convert the cell to a Number. Force the
# cell to Text format first, write the text, then restore the default
# "Normal" style so no visible formatting change is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D5").Value = "599.23"
$ws.Range("D6").Value = "175.17"
$ws.Range("D17").Value = "26.28"
$ws.Range("D19").Value = "12.27"
$ws.Range("D20").Value = "8.17"
$ws.Range("D21").Value = "372.10"
$ws.Range("D24").Value = "72.34"
$ws.Range("D27").Value = "9.79"
$ws.Range("D32").Value = "503.66"
$ws.Range("D34").Value = "1.83"
$ws.Range("D36").Value = "163.44"
$ws.Range("D37").Value = "19.65"
$ws.Range("D46").Value = "156.99"
$ws.Range("D47").Value = "39.51"
$ws.Range("D48").Value = "0.566"
$ws.Range("D50").Value = "1.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
